$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns.Item(4).Insert()
Write-Host "After insert, D7 value:" $ws.Range("D7").Value2
Write-Host "After insert, E7 value:" $ws.Range("E7").Value2
Write-Host "After insert, D7 numfmt:" $ws.Range("D7").NumberFormat
Write-Host "After insert, D8 numfmt:" $ws.Range("D8").NumberFormat
Write-Host "After insert, D102 numfmt:" $ws.Range("D102").NumberFormat
